# "Cmbios datos excel, Se agregaron capturas, Se monitoriaron errores en
#  gestion Facturacion, otros"
#
# The "montoSecundario" (column I) captures for the FACTURA ELECTRONICA
# EXENTA / DISTRIBUCION rows were corrected: the old captured amounts
# ("690,56" for Compra rows, "695,56" for Venta rows) are replaced with
# the newly monitored value "640,56" for rows 2, 3, 6 and 7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "640,56"
$ws.Range("I3").Value = "640,56"
$ws.Range("I6").Value = "640,56"
$ws.Range("I7").Value = "640,56"

# Leave the selection on the cells that were just reviewed/edited (I3, I7,
# I6), ending with I6 as the active cell - matching where the user's
# review/edit pass finished.
$ws.Range("I3").Select()
$ws.Range("I7").Select()
$ws.Range("I6").Select()
